# Visual Performance.xlsx — "Updated visual chart (new test with more controls)"
#
# Updates the two "Milliseconds" sample values (Without Visual / Visual) in the
# results table, marks the "MegaBytes" value for the "Without Visual" row (D4)
# as underlined (new test annotation), and leaves the selection on D11 (where
# the next batch of "more controls" data will be entered).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New measured values for the Milliseconds column.
$ws.Range("C4").Value = 612
$ws.Range("C5").Value = 625

# Flag the MegaBytes figure for the "Without Visual" row with an underline
# (highlighting it as part of the new test run with more controls).
$ws.Range("D4").Font.Underline = $true

# Leave the cursor on D11, ready for the next entries.
$ws.Range("D11").Select()
